$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Reorder sheets: move "M-1 Tasks" to be the first sheet in the workbook.
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("M-1 Tasks").Move($wb.Worksheets.Item(1))

# ---------------------------------------------------------------------------
# 2) "M1 - Game Data" sheet content edits.
# ---------------------------------------------------------------------------
$wsGame = $wb.Worksheets.Item("M1 - Game Data")

# New admin-feature rows (2 & 3 were blank separator rows under the header).
$wsGame.Range("A2").Value = "Admin can edit a game entry"
$wsGame.Range("C2").Value = "Done"
$wsGame.Range("C2").Style = "Good"

$wsGame.Range("A3").Value = "Admin can delete a game entry"
$wsGame.Range("C3").Value = "Done"
$wsGame.Range("C3").Style = "Good"

# "User can enter game data" moves from "In progress" to "Done".
$wsGame.Range("C9").Value = "Done"
$wsGame.Range("C9").Style = "Good"

# Insert two new user-feature rows right after the "enter game data" row.
$wsGame.Rows("10:11").Insert()
$wsGame.Range("C10:C11").Clear()

$wsGame.Range("A10").Value = "User can edit data for their own game."
$wsGame.Range("D10").Value = "Need to restrict this to the owning user"

$wsGame.Range("A11").Value = "User can delete their own game entry"
$wsGame.Range("D11").Value = "Need to restrict this to the owning user"

# Note added to the "filtered by game" row (now pushed down to row 13).
$wsGame.Range("D13").Value = "Change to View Model that shows names. Also need to add filter. (Route? )"

# "View details for a single game entry" row (now row 16) is marked Done.
$wsGame.Range("C16").Value = "Done"
$wsGame.Range("C16").Style = "Good"

# Two new MISC rows appended at the bottom of the sheet.
$wsGame.Range("A24").Value = "Ability to add army list field to WarmachineGame model"
$wsGame.Range("A25").Value = "Expand/Collapse non-required fields in the entry form."

# ---------------------------------------------------------------------------
# 3) View-state updates (selection on "M0 - Account Mgmt", then make
#    "M1 - Game Data" the active/selected tab).
# ---------------------------------------------------------------------------
$wsAccount = $wb.Worksheets.Item("M0 - Account Mgmt")
$wsAccount.Activate()
$wsAccount.Range("D31").Select()

$wsGame.Activate()
$wsGame.Range("A29").Select()
